$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "Goebacillus thermoglucosidasius M10EXG"
$ws.Range("B2").Value = 2501416905
$ws.Range("C2").Value = "Facultative anaerobe"

# Delete rows 3 through 8
$ws.Range("A3:D8").EntireRow.Delete()

# Update selection
$ws.Range("B2").Select()
